# Docx writer: Use different style for block quotes in notes.
# Add a "Footnote Block Text" paragraph style, based on "Footnote Text",
# mirroring the existing "Block Text" style (which is based on "Body Text").

$d = $word.ActiveDocument

$s = $d.Styles.Add("Footnote Block Text", 1)

$s.BaseStyle = "Footnote Text"
$s.NextParagraphStyle = "Footnote Text"
$s.Priority = 9
$s.UnhideWhenUsed = $true
$s.QuickStyle = $true

$pf = $s.ParagraphFormat
$pf.SpaceBefore = 5
$pf.SpaceAfter = 5
$pf.FirstLineIndent = 0
$pf.LeftIndent = 24
$pf.RightIndent = 24
